# Updated symbol list on Mon Dec 26 11:40:21 UTC 2022 with GitHub Actions
#
# The "Price" column (D) and "Volume(1h)" column (E) hold text values
# (e.g. "243.50", "0.006184", "40KickTokenKICK"). Numeric-looking cells
# need to be forced to Text format first (NumberFormat "@") so Excel
# doesn't silently convert the assigned string into a real number and
# drop significant/trailing zeros.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - BNB
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "243.49"

# Row 4 - HuobiToken
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.403"

# Row 5 - Cronos
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.05912"

# Row 6 - GateToken
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "3.456"

# Row 7 - KuCoinToken
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "6.557"

# Row 9 - FTXToken
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9151"

# Row 10 - WazirX
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1418"

# Row 11 - MandalaExchangeToken
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07388"

# Row 12 - LiechtensteinCryptoassetsExchange
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03272"

# Row 14 - BitMartToken
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09353"

# Row 15 - MCDex
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.857"

# Row 16 - BitForexToken
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001574"

# Row 18 - One (drop "Worstin24h" suffix from volume label)
$ws.Range("E18").Value = "17OneONE"

# Row 19 - TigerCash
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.005952"

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0009809"

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.00008604"

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.3240"

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0002902"

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03964"

# Row 41 - KickToken (add "Bestin24h" suffix to volume label)
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006190"
$ws.Range("E41").Value = "40KickTokenKICKBestin24h"

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1075"

# Row 44 - LocalTraders (drop "Bestin24h" suffix from volume label)
$ws.Range("E44").Value = "43LocalTradersLCT"

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005179"

# Row 47 - CoinbaseStockToken (add "Worstin24h" suffix to volume label)
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.8795"
$ws.Range("E47").Value = "46CoinbaseStockTokenCOINWorstin24h"

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.002266"

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00002101"

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0002001"
